# Updates the Hangzhou comic-con workbook ("杭州-漫展信息.xlsx") to the
# refreshed data snapshot:
#   - Sheet "展览": bump several view-count (column F) values, replace the
#     event previously in row 46 ("杭州·梦漫星河动漫展") with a brand-new
#     event ("杭州·AP动漫游戏嘉年华"), and re-insert the original event as a
#     new row 47 (pushing the old row 47 down to row 48, with its own F
#     value bumped too).
#   - Sheet "演出": bump two view-count (column F) values.
#   - Sheet "全部类型": bump the same view-count values as above (this sheet
#     is an independent union/export of all events, refreshed separately;
#     it neither gains a row nor has any text changed).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------
# 1. Sheet "展览" — column F (view count) refreshes for existing rows.
# ---------------------------------------------------------------------
$exhibitFUpdates = @{
    3  = 459
    4  = 211
    7  = 1302
    10 = 343
    11 = 174
    12 = 213
    15 = 16
    17 = 88
    18 = 239
    19 = 1647
    20 = 608
    22 = 189
    23 = 1848
    24 = 399
    26 = 917
    27 = 1202
    30 = 2807
    31 = 1595
    32 = 80
    34 = 629
    35 = 856
    36 = 1756
    37 = 882
    38 = 1775
    39 = 197
    41 = 834
    42 = 36
    43 = 832
    44 = 784
}
foreach ($row in $exhibitFUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitFUpdates[$row]
}

# ---------------------------------------------------------------------
# 2. Sheet "展览" — insert a new row 47 for the returning event, shifting
#    the old row 47 ("理想乡动漫展") down to row 48.
# ---------------------------------------------------------------------
$wsExhibit.Rows.Item(47).Insert()

# Row 47 is blank after the insert; give column A the same look
# (bold / bordered / centered) as the rest of that column.
$wsExhibit.Range("A46").Copy() | Out-Null
$wsExhibit.Range("A47").PasteSpecial(-4122) | Out-Null

# New row 47 gets the data that used to live in row 46
# ("杭州·梦漫星河动漫展").
$wsExhibit.Range("A47").Value = 46
# B47 holds a literal date-like string ("2024-08-03"), not a real date —
# force text parsing (NumberFormat "@") and then restore the cell to the
# default "Normal" style afterwards so no stray number-format sticks
# around on the saved cell (matches the rest of the date column).
$wsExhibit.Range("B47").NumberFormat = "@"
$wsExhibit.Range("B47").Value = "2024-08-03"
$wsExhibit.Range("B47").Style = "Normal"
$wsExhibit.Range("C47").Value = "杭州·梦漫星河动漫展"
$wsExhibit.Range("D47").Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$wsExhibit.Range("E47").Value = "2024.08.03 10:00-08.04 17:00"
$wsExhibit.Range("F47").Value = 428
$wsExhibit.Range("G47").Value = 68
$wsExhibit.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=82836"
$wsExhibit.Range("I47").Value = "//i0.hdslb.com/bfs/openplatform/202403/VFfQUJdD1711700169290.jpeg"

# ---------------------------------------------------------------------
# 3. Sheet "展览" — row 46 now becomes the new event
#    ("杭州·AP动漫游戏嘉年华").
# ---------------------------------------------------------------------
$wsExhibit.Range("C46").Value = "杭州·AP动漫游戏嘉年华"
$wsExhibit.Range("D46").Value = "沈半路171号 Tcar汽车文化主题公园"
$wsExhibit.Range("E46").Value = "2024.08.03 09:00-08.04 17:00"
$wsExhibit.Range("F46").Value = 0
$wsExhibit.Range("G46").Value = 70
$wsExhibit.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=85527"
$wsExhibit.Range("I46").Value = "//i0.hdslb.com/bfs/openplatform/202405/aksNVlQ11715223010331.jpeg"

# ---------------------------------------------------------------------
# 4. Sheet "展览" — row 48 (the shifted-down "理想乡动漫展" row) keeps its
#    own text/values but its row index (col A) and view count (col F)
#    need to reflect the new position / refreshed count.
# ---------------------------------------------------------------------
$wsExhibit.Range("A48").Value = 47
$wsExhibit.Range("F48").Value = 3315

# ---------------------------------------------------------------------
# 5. Sheet "演出" — column F (view count) refreshes.
# ---------------------------------------------------------------------
$wsShow.Cells.Item(3, 6).Value = 186
$wsShow.Cells.Item(12, 6).Value = 788

# ---------------------------------------------------------------------
# 6. Sheet "全部类型" — column F (view count) refreshes (independent
#    snapshot of the same underlying events; no row insertion here).
# ---------------------------------------------------------------------
$allFUpdates = @{
    3  = 459
    4  = 211
    7  = 186
    8  = 1302
    11 = 346
    12 = 174
    13 = 213
    16 = 16
    18 = 88
    20 = 239
    21 = 1647
    22 = 608
    24 = 189
    25 = 1848
    26 = 399
    28 = 1202
    29 = 2807
    30 = 1595
    31 = 80
    33 = 788
    35 = 629
    36 = 856
    37 = 1756
    39 = 882
    40 = 1775
    41 = 834
    42 = 832
    43 = 784
    48 = 3315
}
foreach ($row in $allFUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allFUpdates[$row]
}
